# Update column F (dSF) values on the active sheet per repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -8
$ws.Range("F3").Value = -1
$ws.Range("F5").Value = -9
$ws.Range("F6").Value = -6
$ws.Range("F7").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("F26").Value = -10
$ws.Range("F27").Value = -6
$ws.Range("F31").Value = -4
$ws.Range("F37").Value = -6
